$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.00223048327137546
$ws.Range("D2").Value = 0.0423791821561338
$ws.Range("E2").Value = 0.00223048327137546
$ws.Range("F2").Value = 0.992565055762082
$ws.Range("G2").Value = 0.041635687732342
$ws.Range("H2").Value = 0.00371747211895911
$ws.Range("I2").Value = 0.888475836431227
$ws.Range("L2").Value = 0.642379182156134
$ws.Range("M2").Value = 0.00223048327137546
$ws.Range("N2").Value = 0.980669144981413
$ws.Range("O2").Value = 0.985873605947955
$ws.Range("P2").Value = 0.179182156133829
$ws.Range("Q2").Value = 0.568029739776952
$ws.Range("R2").Value = 0.990334572490706
$ws.Range("S2").Value = 0.00297397769516729
$ws.Range("T2").Value = 0.00297397769516729
$ws.Range("U2").Value = 0.0111524163568773
$ws.Range("V2").Value = 0.150185873605948
$ws.Range("W2").Value = 0.00520446096654275
$ws.Range("X2").Value = 0.0111524163568773

# Row 3
$ws.Range("B3").Value = 0.993308550185874
$ws.Range("C3").Value = 0.89814126394052
$ws.Range("D3").Value = 0.00520446096654275
$ws.Range("E3").Value = 0.00594795539033457
$ws.Range("F3").Value = 0.000743494423791822
$ws.Range("G3").Value = 0.0148698884758364
$ws.Range("J3").Value = 0.987360594795539
$ws.Range("K3").Value = 0.00148698884758364
$ws.Range("L3").Value = 0.000743494423791822
$ws.Range("M3").Value = 0.950185873605948
$ws.Range("N3").Value = 0.00446096654275093
$ws.Range("O3").Value = 0.00446096654275093
$ws.Range("P3").Value = 0.00817843866171004
$ws.Range("Q3").Value = 0.00223048327137546
$ws.Range("R3").Value = 0.00520446096654275
$ws.Range("S3").Value = 0.863940520446097
$ws.Range("T3").Value = 0.00817843866171004
$ws.Range("V3").Value = 0.739776951672863
$ws.Range("W3").Value = 0.00297397769516729
$ws.Range("X3").Value = 0.0163568773234201

# Row 4
$ws.Range("B4").Value = 0.00148698884758364
$ws.Range("C4").Value = 0.00594795539033457
$ws.Range("D4").Value = 0.850557620817844
$ws.Range("E4").Value = 0.991078066914498
$ws.Range("F4").Value = 0.00297397769516729
$ws.Range("G4").Value = 0.937546468401487
$ws.Range("H4").Value = 0.00297397769516729
$ws.Range("I4").Value = 0.0921933085501859
$ws.Range("J4").Value = 0.000743494423791822
$ws.Range("L4").Value = 0.00594795539033457
$ws.Range("M4").Value = 0.00371747211895911
$ws.Range("N4").Value = 0.0118959107806691
$ws.Range("O4").Value = 0.00446096654275093
$ws.Range("P4").Value = 0.811895910780669
$ws.Range("Q4").Value = 0.331598513011152
$ws.Range("R4").Value = 0.00446096654275093
$ws.Range("S4").Value = 0.0966542750929368
$ws.Range("T4").Value = 0.00297397769516729
$ws.Range("U4").Value = 0.984386617100372
$ws.Range("V4").Value = 0.00892193308550186
$ws.Range("W4").Value = 0.107063197026022
$ws.Range("X4").Value = 0.883271375464684

# Row 5
$ws.Range("B5").Value = 0.00520446096654275
$ws.Range("C5").Value = 0.0936802973977695
$ws.Range("D5").Value = 0.10185873605948
$ws.Range("F5").Value = 0.00371747211895911
$ws.Range("G5").Value = 0.00594795539033457
$ws.Range("H5").Value = 0.993308550185874
$ws.Range("I5").Value = 0.0193308550185874
$ws.Range("J5").Value = 0.0118959107806691
$ws.Range("K5").Value = 0.998513011152416
$ws.Range("L5").Value = 0.35092936802974
$ws.Range("M5").Value = 0.0438661710037175
$ws.Range("N5").Value = 0.00297397769516729
$ws.Range("O5").Value = 0.00520446096654275
$ws.Range("P5").Value = 0.000743494423791822
$ws.Range("Q5").Value = 0.0981412639405204
$ws.Range("S5").Value = 0.0356877323420074
$ws.Range("T5").Value = 0.985873605947955
$ws.Range("U5").Value = 0.00446096654275093
$ws.Range("V5").Value = 0.101115241635688
$ws.Range("W5").Value = 0.884758364312268
$ws.Range("X5").Value = 0.0892193308550186
